$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L: "E" (energy) header plus a constant value of 252 for every
# data row (rows 2-149), mirroring the existing K column's per-row constant.
$ws.Range("L1").Value = "E"
$ws.Range("L2:L149").Value = 252

# Move the active selection to the newly populated column, matching the
# saved workbook state after the edit.
[void]$ws.Range("L2:L149").Select()
